$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sales rows appended beneath the existing data (row 2 is the last
# populated row). Values that look numeric ("121210", "800", "4", "3200",
# "2400") must be stored as literal text (shared strings), matching every
# other data cell in this sheet - so each row is pre-formatted as Text,
# filled in, and then restored to the Normal style so no stray formatting
# is left behind.

$ws.Range("A3:F4").NumberFormat = "@"

$ws.Range("A3").Value = "121210"
$ws.Range("B3").Value = "2023-10-31"
$ws.Range("C3").Value = "Oral B"
$ws.Range("D3").Value = "800"
$ws.Range("E3").Value = "4"
$ws.Range("F3").Value = "3200"

$ws.Range("A4").Value = "121210"
$ws.Range("B4").Value = "2023-10-31"
$ws.Range("C4").Value = "Oral B"
$ws.Range("D4").Value = "800"
$ws.Range("E4").Value = "3"
$ws.Range("F4").Value = "2400"

$ws.Range("A3:F4").Style = "Normal"
